$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos feed refresh: Price (D) and Volume(1h) (E) columns updated,
# plus a few re-ranked coins (rows 41-44) where Coin/Link/Price/Volume all moved.
# Price values are plain numeric-looking strings (e.g. "1.00", "525.73",
# "21.07") that Excel would otherwise auto-convert to real numbers on
# assignment, so those cells are pre-formatted as Text to keep the exact
# literal string (matching thousand-dot-separated values like "57.882.48"
# which already round-trip fine as text).

$ws.Range('D2').Value = '57.882.48'
$ws.Range('E2').Value = '  +2.40%  '
$ws.Range('D3').Value = '3.055.23'
$ws.Range('E3').Value = '  +2.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.73'
$ws.Range('E5').Value = '  +5.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.75'
$ws.Range('E6').Value = '  +5.54%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +4.90%  '
$ws.Range('E9').Value = '  +5.20%  '
$ws.Range('E10').Value = '  +7.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.367'
$ws.Range('E11').Value = '  +4.83%  '
$ws.Range('E12').Value = '  +2.37%  '
$ws.Range('D13').Value = '3.576.93'
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.40'
$ws.Range('E14').Value = '  +6.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000168'
$ws.Range('E15').Value = '  +15.23%  '
$ws.Range('D16').Value = '57.908.51'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.16'
$ws.Range('E17').Value = '  +6.30%  '
$ws.Range('D18').Value = '3.051.04'
$ws.Range('E18').Value = '  +2.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.97'
$ws.Range('E19').Value = '  +6.32%  '
$ws.Range('E20').Value = '  +6.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '340.23'
$ws.Range('E21').Value = '  +5.58%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.498'
$ws.Range('E23').Value = '  +8.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.36'
$ws.Range('E24').Value = '  +6.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.174'
$ws.Range('E25').Value = '  +7.07%  '
$ws.Range('D26').Value = '0.0₃0966'
$ws.Range('E26').Value = '  +8.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.94'
$ws.Range('E28').Value = '  +6.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.27'
$ws.Range('E29').Value = '  +7.76%  '
$ws.Range('E30').Value = '  +7.71%  '
$ws.Range('E31').Value = '  +5.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.07'
$ws.Range('E32').Value = '  +6.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.26'
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.73'
$ws.Range('E34').Value = '  +6.25%  '
$ws.Range('E35').Value = '  +5.77%  '
$ws.Range('E36').Value = '  +3.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.92'
$ws.Range('E37').Value = '  +11.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0692'
$ws.Range('E38').Value = '  +4.02%  '
$ws.Range('D39').Value = '3.087.37'
$ws.Range('E39').Value = '  +2.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.65'
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.84'
$ws.Range('E42').Value = '  +8.25%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.343.64'
$ws.Range('E43').Value = '  +6.20%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.48'
$ws.Range('E44').Value = '  +4.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.661'
$ws.Range('E45').Value = '  +4.45%  '
$ws.Range('E46').Value = '  +3.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.04'
$ws.Range('E47').Value = '  +4.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.03'
$ws.Range('E49').Value = '  +7.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.08'
$ws.Range('E50').Value = '  +5.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0893'
$ws.Range('E51').Value = '  +5.46%  '
